# Applies the row-data swaps described in the diff:
#   - Row 10  <-> Row 12   (species record swap)
#   - Row 13  <-> Row 14   (species record swap)
#   - Row 27  <-> Row 28   (species record swap)
# Each block below writes the full target row values (columns
# A,B,D,E,F,G,H,I,J,Q,R,X,Z,AB) taken from the unified diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 10 becomes what used to be row 12's data ----
# (Rosenticka / Fomitopsis rosea, Antal=2)
$ws.Range("A10").Value = 131106319
$ws.Range("B10").Value = 92107
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 658
$ws.Range("F10").Value = "Rosenticka"
$ws.Range("G10").Value = "Fomitopsis rosea"
$ws.Range("H10").Value = "(Alb. & Schwein.:Fr.) P.Karst."
$ws.Range("I10").Value = "'2"
$ws.Range("Q10").Value = 601569
$ws.Range("R10").Value = 6992657
$ws.Range("X10").Value = "2025_0864"
$ws.Range("Z10").Value = "13:14"
$ws.Range("AB10").Value = "13:14"

# ---- Row 12 becomes what used to be row 10's data ----
# (Lappticka / Amylocystis lapponica)
$ws.Range("A12").Value = 131106321
$ws.Range("B12").Value = 92022
$ws.Range("D12").Value = "VU"
$ws.Range("E12").Value = 48
$ws.Range("F12").Value = "Lappticka"
$ws.Range("G12").Value = "Amylocystis lapponica"
$ws.Range("H12").Value = "(Romell) Bondartsev & Singer"
$ws.Range("I12").Value = ""
$ws.Range("Q12").Value = 601579
$ws.Range("R12").Value = 6992698
$ws.Range("X12").Value = "2025_0862"
$ws.Range("Z12").Value = "12:56"
$ws.Range("AB12").Value = "12:56"

# ---- Row 13 becomes what used to be row 14's data ----
# (Rosenticka / Fomitopsis rosea)
$ws.Range("A13").Value = 131106312
$ws.Range("B13").Value = 92107
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 658
$ws.Range("F13").Value = "Rosenticka"
$ws.Range("G13").Value = "Fomitopsis rosea"
$ws.Range("H13").Value = "(Alb. & Schwein.:Fr.) P.Karst."
$ws.Range("I13").Value = ""
$ws.Range("Q13").Value = 601540
$ws.Range("R13").Value = 6992576
$ws.Range("X13").Value = "2025_0872"
$ws.Range("Z13").Value = "13:29"
$ws.Range("AB13").Value = "13:29"

# ---- Row 14 becomes what used to be row 13's data ----
# (Ullticka / Phellinidium ferrugineofuscum)
$ws.Range("A14").Value = 131106325
$ws.Range("B14").Value = 91809
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 1202
$ws.Range("F14").Value = "Ullticka"
$ws.Range("G14").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H14").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("I14").Value = ""
$ws.Range("Q14").Value = 601615
$ws.Range("R14").Value = 6992785
$ws.Range("X14").Value = "2025_0858"
$ws.Range("Z14").Value = "12:21"
$ws.Range("AB14").Value = "12:21"

# ---- Row 27 becomes what used to be row 28's data ----
# (Rosenticka / Fomitopsis rosea); loses its "mycel" J value
$ws.Range("A27").Value = 131106330
$ws.Range("B27").Value = 92107
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 658
$ws.Range("F27").Value = "Rosenticka"
$ws.Range("G27").Value = "Fomitopsis rosea"
$ws.Range("H27").Value = "(Alb. & Schwein.:Fr.) P.Karst."
$ws.Range("I27").Value = ""
$ws.Range("J27").Value = ""
$ws.Range("Q27").Value = 601607
$ws.Range("R27").Value = 6992782
$ws.Range("X27").Value = "2025_0853"
$ws.Range("Z27").Value = "12:06"
$ws.Range("AB27").Value = "12:06"

# ---- Row 28 becomes what used to be row 27's data ----
# (Rynkskinn / Hermanssonia centrifuga); gains "mycel" J value
$ws.Range("A28").Value = 131106329
$ws.Range("B28").Value = 92268
$ws.Range("D28").Value = "VU"
$ws.Range("E28").Value = 1209
$ws.Range("F28").Value = "Rynkskinn"
$ws.Range("G28").Value = "Hermanssonia centrifuga"
$ws.Range("H28").Value = "(P. Karst.) Zmitr."
$ws.Range("I28").Value = ""
$ws.Range("J28").Value = "mycel"
$ws.Range("Q28").Value = 601609
$ws.Range("R28").Value = 6992789
$ws.Range("X28").Value = "2025_0854"
$ws.Range("Z28").Value = "12:09"
$ws.Range("AB28").Value = "12:09"

Write-Host "Row swaps applied."
